$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 9.569019626703195
    "C2" = 4.77145487767217
    "D2" = 9.111797131698783
    "E2" = 13.73263851607108
    "F2" = 33.53764191898373
    "I2" = 23.12770959703998
    "J2" = 10.16022615965016
    "K2" = 9.879565221039705
    "M2" = 15.27455908910457
    "N2" = 19.72310094328491
    "O2" = 25.33748623929682
    "B3" = 9.323129051461319
    "C3" = 4.608940374867309
    "D3" = 9.081071464213434
    "E3" = 13.72823215512705
    "F3" = 33.59627219694943
    "I3" = 23.20840598797397
    "J3" = 10.17970882645873
    "K3" = 9.72431153211933
    "M3" = 15.21705972128952
    "N3" = 19.78037058264195
    "O3" = 25.40848426633846
    "B4" = 9.170404399642626
    "C4" = 4.506845230006216
    "D4" = 9.063625518905512
    "E4" = 13.72782019895599
    "F4" = 33.63918473314339
    "I4" = 23.26173951049202
    "J4" = 10.19279432961303
    "K4" = 9.629292081229572
    "M4" = 15.18400478770095
    "N4" = 19.81719153603212
    "O4" = 25.45667817026714
    "B5" = 9.107823740127978
    "C5" = 4.464724283594077
    "D5" = 9.0568784276144
    "E5" = 13.72823079897826
    "F5" = 33.65840857046877
    "I5" = 23.28442515574811
    "J5" = 10.19840949986023
    "K5" = 9.590695588802371
    "M5" = 15.17111028186028
    "N5" = 19.83261428344701
    "O5" = 25.47747334745744
    "B6" = 9.09741435629485
    "C6" = 4.457700947733464
    "D6" = 9.055780105576938
    "E6" = 13.72833396301123
    "F6" = 33.66170550223953
    "I6" = 23.28824957304356
    "J6" = 10.19935897946649
    "K6" = 9.584295596422464
    "M6" = 15.16900420395823
    "N6" = 19.83520049503994
    "O6" = 25.48099614154455
    "B7" = 9.169561679193622
    "C7" = 4.506279173232492
    "D7" = 9.063533051743342
    "E7" = 13.72782339211653
    "F7" = 33.63943696364459
    "I7" = 23.26204160338371
    "J7" = 10.19286891261677
    "K7" = 9.628770987526766
    "M7" = 15.18382854448871
    "N7" = 19.81739783890121
    "O7" = 25.45695394320491
    "B8" = 9.48465470477362
    "C8" = 4.715936477412281
    "D8" = 9.100911520049872
    "E8" = 13.73064429063252
    "F8" = 33.55642183003316
    "I8" = 23.15474793697524
    "J8" = 10.16671088300262
    "K8" = 9.825998445407086
    "M8" = 15.25427244606797
    "N8" = 19.74250429799471
    "O8" = 25.36101079040478
    "B9" = 10.08476391724651
    "C9" = 5.10618764783627
    "D9" = 9.185229100752917
    "E9" = 13.75428849249351
    "F9" = 33.44854839740529
    "I9" = 22.97438779415082
    "J9" = 10.12431202030651
    "K9" = 10.21313256478079
    "M9" = 15.40980522902307
    "N9" = 19.60873130666223
    "O9" = 25.20942116238844
    "B10" = 10.5099479312334
    "C10" = 5.377148936399037
    "D10" = 9.253543277860073
    "E10" = 13.78257224336031
    "F10" = 33.40284094134115
    "I10" = 22.86019884663709
    "J10" = 10.09856623543414
    "K10" = 10.4950804859547
    "M10" = 15.53402794524846
    "N10" = 19.51835064114272
    "O10" = 25.12039676418447
    "B11" = 10.69904819112883
    "C11" = 5.496463380679391
    "D11" = 9.285922184016634
    "E11" = 13.79777566831007
    "F11" = 33.3893369596964
    "I11" = 22.81223149901389
    "J11" = 10.08802319577616
    "K11" = 10.62225782871494
    "M11" = 15.59255218842403
    "N11" = 19.47893309572052
    "O11" = 25.08476207624785
    "B12" = 10.76996648470068
    "C12" = 5.54103942822212
    "D12" = 9.298363194873428
    "E12" = 13.80386558082717
    "F12" = 33.38527095118873
    "I12" = 22.79463976348449
    "J12" = 10.08419855131266
    "K12" = 10.67021719756004
    "M12" = 15.6149901929657
    "N12" = 19.46424948458625
    "O12" = 25.07196828513288
    "B13" = 10.75472467813321
    "C13" = 5.531466678384484
    "D13" = 9.295675924594452
    "E13" = 13.80253926835007
    "F13" = 33.38610005646023
    "I13" = 22.79840299089078
    "J13" = 10.08501479997301
    "K13" = 10.65989786709596
    "M13" = 15.61014570590878
    "N13" = 19.46740107442287
    "O13" = 25.07469249643902
    "B14" = 10.7048968816483
    "C14" = 5.500143021202056
    "D14" = 9.286942153257508
    "E14" = 13.79827004134592
    "F14" = 33.38898145444919
    "I14" = 22.81077274078321
    "J14" = 10.08770517926302
    "K14" = 10.62620772897552
    "M14" = 15.59439272136019
    "N14" = 19.47772020310208
    "O14" = 25.08369548308896
    "B15" = 10.67428413135409
    "C15" = 5.48087643714875
    "D15" = 9.28161565820286
    "E15" = 13.79569824040531
    "F15" = 33.39088280805464
    "I15" = 22.81842413441623
    "J15" = 10.08937495448724
    "K15" = 10.60554422617459
    "M15" = 15.58477911602495
    "N15" = 19.48407257581284
    "O15" = 25.08930129902418
    "B16" = 10.4974968650399
    "C16" = 5.369268663383704
    "D16" = 9.251452816277371
    "E16" = 13.78162540996276
    "F16" = 33.40387009788902
    "I16" = 22.86341370768725
    "J16" = 10.09927874236584
    "K16" = 10.48674354130732
    "M16" = 15.53024263877545
    "N16" = 19.52096072738775
    "O16" = 25.12282351827098
    "B17" = 10.387886699055
    "C17" = 5.299761322298505
    "D17" = 9.233277247337078
    "E17" = 13.77358844924811
    "F17" = 33.41370401388614
    "I17" = 22.89203248174961
    "J17" = 10.10565355243575
    "K17" = 10.41355389128384
    "M17" = 15.49729276797785
    "N17" = 19.54402430340318
    "O17" = 25.14463458227162
    "B18" = 10.32443913247786
    "C18" = 5.259412982566952
    "D18" = 9.222946197332281
    "E18" = 13.76918599020663
    "F18" = 33.42004631514926
    "I18" = 22.90886762295354
    "J18" = 10.10943021065525
    "K18" = 10.37135782060444
    "M18" = 15.47853134674138
    "N18" = 19.55744966264984
    "O18" = 25.15763740044116
    "B19" = 10.3028899108737
    "C19" = 5.245689494512256
    "D19" = 9.219469639626892
    "E19" = 13.76773330839852
    "F19" = 33.42231155119806
    "I19" = 22.91463199866701
    "J19" = 10.11072783012594
    "K19" = 10.35705527795623
    "M19" = 15.47221217250785
    "N19" = 19.56202273753596
    "O19" = 25.16211850437184
    "B20" = 10.39959708430398
    "C20" = 5.307199035877548
    "D20" = 9.235199385066815
    "E20" = 13.77442123220492
    "F20" = 33.41258617017828
    "I20" = 22.88894721423686
    "J20" = 10.10496355686764
    "K20" = 10.4213556605336
    "M20" = 15.50078072173983
    "N20" = 19.54155261547343
    "O20" = 25.1422653824173
    "B21" = 10.71955175621035
    "C21" = 5.509360259677662
    "D21" = 9.289502654084432
    "E21" = 13.79951501535732
    "F21" = 33.38810669092837
    "I21" = 22.8071239053588
    "J21" = 10.08691039871995
    "K21" = 10.63610910920007
    "M21" = 15.59901237042115
    "N21" = 19.4746826396111
    "O21" = 25.08103207448022
    "B22" = 10.92461170867865
    "C22" = 5.637939019957971
    "D22" = 9.32603747905714
    "E22" = 13.81785276707416
    "F22" = 33.37821391089078
    "I22" = 22.75698458370344
    "J22" = 10.07608939538163
    "K22" = 10.77527834090983
    "M22" = 15.66481544439961
    "N22" = 19.43239489575251
    "O22" = 25.04509479418772
    "B23" = 10.81555907923559
    "C23" = 5.569649949150372
    "D23" = 9.30644511578573
    "E23" = 13.8078894596264
    "F23" = 33.38293545343491
    "I23" = 22.78343939975617
    "J23" = 10.08177540404343
    "K23" = 10.70112364271042
    "M23" = 15.62955301259025
    "N23" = 19.45483548215886
    "O23" = 25.0639013690678
    "B24" = 10.39430415698709
    "C24" = 5.303837650046004
    "D24" = 9.2343300174705
    "E24" = 13.77404405144003
    "F24" = 33.41308940219276
    "I24" = 22.89034087434305
    "J24" = 10.10527515592265
    "K24" = 10.41782884594057
    "M24" = 15.499203250224
    "N24" = 19.54266954903281
    "O24" = 25.14333505444775
    "B25" = 9.924848208738704
    "C25" = 5.003183817868613
    "D25" = 9.161274303971169
    "E25" = 13.74596532190499
    "F25" = 33.47184279514008
    "I25" = 23.0199630645842
    "J25" = 10.13483150711945
    "K25" = 10.10864642999077
    "M25" = 15.36593250398706
    "N25" = 19.64352700501365
    "O25" = 25.24650981154413
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

Write-Host "Updated $($values.Count) cells"
